$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we are about to write so Excel
# does not silently reinterpret numeric-looking strings as numbers
# (these columns store text values like "28.029.69" / "  +0.28%  ").
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.029.69"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.864.35"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "312.13"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "0.5093"
$ws.Range("E7").Value = "  +1.25%  "
$ws.Range("D8").Value = "0.3826"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "0.08310"
$ws.Range("E9").Value = "  -6.88%  "
$ws.Range("D10").Value = "1.115"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "6.221"
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("D13").Value = "20.59"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").Value = "1.859.27"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").Value = "7.210"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "0.00001097"
$ws.Range("D18").Value = "90.86"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "17.73"
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "6.038"
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("D23").Value = "28.037.47"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "11.07"
$ws.Range("E24").Value = "  -3.61%  "
$ws.Range("E25").Value = "  -1.98%  "
$ws.Range("D26").Value = "2.550"
$ws.Range("E26").Value = "  +2.07%  "
$ws.Range("D27").Value = "2.073.59"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").Value = "158.21"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "20.59"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").Value = "125.35"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").Value = "0.1054"
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("D32").Value = "1.039"
$ws.Range("E32").Value = "  -1.54%  "
$ws.Range("D33").Value = "5.847"
$ws.Range("E33").Value = "  +4.36%  "
$ws.Range("D34").Value = "3.598"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").Value = "9.431"
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("D36").Value = "0.06527"
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("D37").Value = "0.02414"
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("D38").Value = "0.2169"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("D39").Value = "1.207"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").Value = "0.6468"
$ws.Range("E40").Value = "  +1.44%  "
$ws.Range("D41").Value = "1.224"
$ws.Range("E41").Value = "  -4.58%  "
$ws.Range("D42").Value = "4.940"
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("D43").Value = "11.23"
$ws.Range("E43").Value = "  -2.34%  "
$ws.Range("D44").Value = "0.6106"
$ws.Range("E44").Value = "  +1.65%  "
$ws.Range("D45").Value = "13.11"
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("D47").Value = "3.671"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").Value = "2.018"
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("D49").Value = "1.209"
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("D50").Value = "120.48"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").Value = "78.48"
$ws.Range("E51").Value = "  -0.70%  "
